# Updated Argent prices: append a new row (2025-04-10) to every price
# sheet, repeating the latest known price for that sheet.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-10"

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "41.5"
    "N-type Wafer"               = "1.28"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,038"
    "Silver Busbar front-side"   = "7,542"
    "Silver finger front-side"   = "7,592"
    "USD_CNY"                    = "7.3769"
}

foreach ($ws in $wb.Worksheets) {
    if (-not $sheetValues.ContainsKey($ws.Name)) {
        continue
    }

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1
    $value = $sheetValues[$ws.Name]

    # Enter the values as formulas that evaluate to string literals, then
    # convert them to plain values in place. This preserves the text data
    # type (so a value like "2025-04-10" or "5,038" is stored exactly as
    # typed, not auto-converted into a date serial number or a plain
    # number) without leaving behind any extra cell formatting/styles.
    $dateCell = $ws.Cells.Item($newRow, 1)
    $dateCell.Formula = "=""" + $newDate + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $valueCell = $ws.Cells.Item($newRow, 2)
    $valueCell.Formula = "=""" + $value + """"
    $valueCell.Copy()
    $valueCell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
